# This script applies a weekly data update to the "Acelga" price sheet.
# A new daily price record is inserted at row 246 (with a new date,
# 2023-02-03 / serial 44960), which pushes all the existing records in
# rows 246-368 down by one row. The record that used to be in row 368
# "falls off the end" and becomes the new row 369.
#
# Only the variable columns (D = Fecha, J = Volumen, K = Precio minimo,
# L = Precio maximo, M = Precio promedio ponderado, P = Precio $/Kg) carry
# row-specific data; all of the other columns (A, B, C, E, F, G, H, I, N,
# O, Q, R) are constant for every data row in this sheet, so they do not
# need to be touched/shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 246
$lastDataRow = 368
$newLastRow = $lastDataRow + 1

# 1) Snapshot the current ("before") values of the columns that shift,
#    for every data row, before we start overwriting anything.
$colD = @{}
$colJ = @{}
$colK = @{}
$colL = @{}
$colM = @{}
$colP = @{}

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $colD[$r] = $ws.Cells.Item($r, 4).Value()
    $colJ[$r] = $ws.Cells.Item($r, 10).Value()
    $colK[$r] = $ws.Cells.Item($r, 11).Value()
    $colL[$r] = $ws.Cells.Item($r, 12).Value()
    $colM[$r] = $ws.Cells.Item($r, 13).Value()
    $colP[$r] = $ws.Cells.Item($r, 16).Value()
}

# 2) Create the new last row (369), copying the constant columns from the
#    previous last row (368) and using the data that "falls off" row 368.
$ws.Range("A" + $newLastRow).Value = $ws.Range("A" + $lastDataRow).Value()
$ws.Range("B" + $newLastRow).Value = $ws.Range("B" + $lastDataRow).Value()
$ws.Range("C" + $newLastRow).Value = $ws.Range("C" + $lastDataRow).Value()
$ws.Range("E" + $newLastRow).Value = $ws.Range("E" + $lastDataRow).Value()
$ws.Range("F" + $newLastRow).Value = $ws.Range("F" + $lastDataRow).Value()
$ws.Range("G" + $newLastRow).Value = $ws.Range("G" + $lastDataRow).Value()
$ws.Range("H" + $newLastRow).Value = $ws.Range("H" + $lastDataRow).Value()
$ws.Range("I" + $newLastRow).Value = $ws.Range("I" + $lastDataRow).Value()
$ws.Range("N" + $newLastRow).Value = $ws.Range("N" + $lastDataRow).Value()
$ws.Range("O" + $newLastRow).Value = $ws.Range("O" + $lastDataRow).Value()
$ws.Range("Q" + $newLastRow).Value = $ws.Range("Q" + $lastDataRow).Value()
$ws.Range("R" + $newLastRow).Value = $ws.Range("R" + $lastDataRow).Value()

$ws.Range("D" + $newLastRow).NumberFormat = $ws.Range("D" + $lastDataRow).NumberFormat()
$ws.Range("D" + $newLastRow).Value = $colD[$lastDataRow]
$ws.Range("J" + $newLastRow).Value = $colJ[$lastDataRow]
$ws.Range("K" + $newLastRow).Value = $colK[$lastDataRow]
$ws.Range("L" + $newLastRow).Value = $colL[$lastDataRow]
$ws.Range("M" + $newLastRow).Value = $colM[$lastDataRow]
$ws.Range("P" + $newLastRow).Value = $colP[$lastDataRow]

# 3) Shift rows lastDataRow .. firstDataRow+1 down by one: row r takes the
#    "before" values that used to belong to row r-1.
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    $prev = $r - 1
    $ws.Cells.Item($r, 4).Value = $colD[$prev]
    $ws.Cells.Item($r, 10).Value = $colJ[$prev]
    $ws.Cells.Item($r, 11).Value = $colK[$prev]
    $ws.Cells.Item($r, 12).Value = $colL[$prev]
    $ws.Cells.Item($r, 13).Value = $colM[$prev]
    $ws.Cells.Item($r, 16).Value = $colP[$prev]
}

# 4) Finally, row firstDataRow (246) receives the brand-new data point.
#    Only its date changes (to 2023-02-03, serial 44960); its other
#    values (J, K, L, M, P) stay the same as they were before.
$ws.Range("D" + $firstDataRow).Value = 44960
